$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily COVID totals rows (342-349), continuing the existing table
# (columns: date, areaType, areaCode, areaName, cumCasesByPublishDate,
#  newCasesByPublishDate, newDeaths28DaysByPublishDate, cumDeaths28DaysByPublishDate)
$data = @(
    @("2021-07-19","overview","K02000001","United Kingdom",5473477,39950,19,128727),
    @("2021-07-20","overview","K02000001","United Kingdom",5519602,46558,96,128823),
    @("2021-07-21","overview","K02000001","United Kingdom",5563006,44104,73,128896),
    @("2021-07-22","overview","K02000001","United Kingdom",5602321,39906,84,128980),
    @("2021-07-23","overview","K02000001","United Kingdom",5637975,36389,64,129044),
    @("2021-07-24","overview","K02000001","United Kingdom",5669260,31795,86,129130),
    @("2021-07-25","overview","K02000001","United Kingdom",5697912,29173,28,129158),
    @("2021-07-26","overview","K02000001","United Kingdom",5722298,24950,14,129172)
)

$startRow = 342

# Force the date column to Text format first so Excel doesn't reinterpret the
# "YYYY-MM-DD" strings as date serial numbers (matches existing column A cells,
# which are stored as literal text).
$ws.Range("A$startRow`:A$($startRow + $data.Count - 1)").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
